$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.830.52"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.889.83"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'0.7788"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "'0.07305"
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").Value = "'25.30"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.08130"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'0.7643"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "'5.456"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").Value = "1.865.04"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "'93.40"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "'6.188"
$ws.Range("E16").Value = "  +4.28%  "
$ws.Range("D17").Value = "29.791.46"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'13.92"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'245.91"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "'0.000007846"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "'0.9993"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'8.138"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.117.55"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "'0.9994"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'0.1586"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("D26").Value = "'9.436"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").Value = "'161.34"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "'18.76"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'2.029"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "'1.448"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("D31").Value = "'1.542"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'4.471"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'0.05584"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "'4.074"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "'0.7530"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("D37").Value = "'0.9963"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'2.637"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").Value = "'2.780"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "1.138.63"
$ws.Range("E41").Value = "  +10.35%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'73.55"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").Value = "'5.948"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").Value = "'0.8567"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'1.897"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").Value = "'101.78"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'3.100"
$ws.Range("E49").Value = "  +6.50%  "
$ws.Range("D50").Value = "'9.773"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "'7.493"
